$d = $word.ActiveDocument

# Add the new paragraph style "Footnote Block Text" (styleId "FootnoteBlockText"),
# based on "Footnote Text" and followed by "Footnote Text", matching the
# "Block Text" style but for use inside footnotes so it can get its own
# font size independent of the surrounding note text.
$s = $d.Styles.Add("FootnoteBlockText", 1)
$s.NameLocal = "Footnote Block Text"
$s.BaseStyle = "Footnote Text"
$s.NextParagraphStyle = "Footnote Text"
$s.Priority = 9
$s.UnhideWhenUsed = $true
$s.QuickStyle = $true

$s.ParagraphFormat.SpaceBefore = 5
$s.ParagraphFormat.SpaceAfter = 5
$s.ParagraphFormat.FirstLineIndent = 0
$s.ParagraphFormat.LeftIndent = 24
$s.ParagraphFormat.RightIndent = 24
